# LH_WF_REGISTRATION_REVIEWS.xlsx - v1.3 verification edit
#
# - close registration wireframe review, verify the updates
#   and modify id naming convention
#
# Changes:
#   Sheet "LH_WF_REGISTRATION_REVIEW":
#     B2: LH_WF_REGISTRATION_REVIEW_001 -> LH-WF-REGISTRATION-REVIEW-001
#     B3: LH_WF_REGISTRATION_REVIEW_002 -> LH-WF-REGISTRATION-REVIEW-002
#     I3: open -> closed   (reviewer verification, now closed)
#   Sheet "VERSION-HISTORY":
#     Row 5 (v1.3 entry):
#       B5: Gehad Ashry -> Ahmed Abuzaid
#       C5: "Closed reviews " -> "close registration wireframe review, verify the updates and modify id naming convention"
#       D5: 29/4/2025 -> 28/4/2025 (45776 -> 45775)

$wb = $excel.ActiveWorkbook

$wsReview  = $wb.Worksheets.Item("LH_WF_REGISTRATION_REVIEW")
$wsHistory = $wb.Worksheets.Item("VERSION-HISTORY")

# --- LH_WF_REGISTRATION_REVIEW sheet: rename IDs to dash convention, close review 002 ---
$wsReview.Range("B2").Value = "LH-WF-REGISTRATION-REVIEW-001"
$wsReview.Range("B3").Value = "LH-WF-REGISTRATION-REVIEW-002"
$wsReview.Range("I3").Value = "closed"

# --- VERSION-HISTORY sheet: update the v1.3 row ---
$wsHistory.Range("B5").Value = "Ahmed Abuzaid"
$wsHistory.Range("C5").Value = "close registration wireframe review, verify the updates and modify id naming convention"
$wsHistory.Range("D5").Value = 45775

# the longer comment text wraps across more lines, so the row grows taller
$wsHistory.Rows.Item(5).RowHeight = 56.25

# --- refresh view selections to match the saved state ---
$wsReview.Range("C3").Select() | Out-Null
$wsHistory.Range("C12:C13").Select() | Out-Null
